$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fase 2")

# --- Fill in RACI letters for the existing activity rows 14-20 ---
$ws.Range("B14").Value = "R"
$ws.Range("C14").Value = "A"
$ws.Range("D14").Value = "C"
$ws.Range("E14").Value = "I"

$ws.Range("B15").Value = "A"
$ws.Range("C15").Value = "I"
$ws.Range("D15").Value = "C"
$ws.Range("E15").Value = "R"

$ws.Range("B16").Value = "I"
$ws.Range("C16").Value = "A"
$ws.Range("D16").Value = "R"
$ws.Range("E16").Value = "C"

$ws.Range("B17").Value = "I"
$ws.Range("C17").Value = "R"
$ws.Range("D17").Value = "A"
$ws.Range("E17").Value = "C"

$ws.Range("B18").Value = "I"
$ws.Range("C18").Value = "A"
$ws.Range("D18").Value = "R"
$ws.Range("E18").Value = "C"

$ws.Range("B19").Value = "I"
$ws.Range("C19").Value = "A"
$ws.Range("D19").Value = "C"
$ws.Range("E19").Value = "R"

$ws.Range("B20").Value = "I"
$ws.Range("C20").Value = "A"
$ws.Range("D20").Value = "C"
$ws.Range("E20").Value = "R"

# --- New activity rows 21-22 ---
$ws.Range("A21").Value = "Apartado Calendario"
$ws.Range("A21").WrapText = $false
$ws.Range("B21").Value = "A"
$ws.Range("C21").Value = "I"
$ws.Range("D21").Value = "R"
$ws.Range("E21").Value = "C"

$ws.Range("A22").Value = "Apartado portal pago"
$ws.Range("A22").WrapText = $false
$ws.Range("B22").Value = "C"
$ws.Range("C22").Value = "R"
$ws.Range("D22").Value = "A"
$ws.Range("E22").Value = "I"
